# Fruta / hortaliza, semanal
# Weekly data update: a new week's record is inserted at the top of the
# "Agricola del Norte S.A. de Arica - Cebollin baby" block (row 58) and
# every subsequent row (59-70) takes on the values that the row above it
# held previously (i.e. the whole block shifts down by one record).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r = row number, values = Fecha(D), Volumen(J), PrecioMinimo(K), PrecioMaximo(L), PrecioPromedio(M), PrecioKg(P)
$rows = @(
    @{ r = 58; D = 44767; J = 250; K = 2500; L = 3000; M = 2750; P = 1375 },
    @{ r = 59; D = 44691; J = 270; K = 2000; L = 2500; M = 2250; P = 1125 },
    @{ r = 60; D = 44529; J = 300; K = 1800; L = 2000; M = 1900; P = 950  },
    @{ r = 61; D = 44533; J = 250; K = 1000; L = 1100; M = 1050; P = 525  },
    @{ r = 62; D = 44687; J = 300; K = 1800; L = 2000; M = 1900; P = 950  },
    @{ r = 63; D = 44627; J = 300; K = 2000; L = 2500; M = 2250; P = 1125 },
    @{ r = 64; D = 44431; J = 300; K = 1900; L = 2000; M = 1950; P = 975  },
    @{ r = 65; D = 44263; J = 270; K = 1900; L = 2000; M = 1950; P = 975  },
    @{ r = 66; D = 44749; J = 300; K = 2000; L = 2500; M = 2250; P = 1125 },
    @{ r = 67; D = 44267; J = 300; K = 1400; L = 1500; M = 1450; P = 725  },
    @{ r = 68; D = 44568; J = 300; K = 5000; L = 5500; M = 5250; P = 2625 },
    @{ r = 69; D = 44736; J = 270; K = 2800; L = 3000; M = 2900; P = 1450 },
    @{ r = 70; D = 44648; J = 300; K = 1300; L = 1500; M = 1400; P = 700  }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("D$r").Value = $row.D
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("P$r").Value = $row.P
}
